$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 2160.906096720676
$ws.Range("C2").Value = 2010.631070832981
$ws.Range("D2").Value = 1347.133603557302
$ws.Range("E2").Value = 2211.062976753502
$ws.Range("F2").Value = 2176.626862922192
$ws.Range("G2").Value = 2136.069371440038
$ws.Range("H2").Value = 2212.565927724553
$ws.Range("B3").Value = 2178.466453297492
$ws.Range("C3").Value = 2018.735400027925
$ws.Range("D3").Value = 1200.691650150008
$ws.Range("E3").Value = 2212.289353747457
$ws.Range("F3").Value = 2186.295744950701
$ws.Range("G3").Value = 2132.944830111488
$ws.Range("H3").Value = 2213.057972484487
$ws.Range("B4").Value = 2122.310463277523
$ws.Range("C4").Value = 2017.413604081494
$ws.Range("D4").Value = 1242.196491198539
$ws.Range("E4").Value = 2204.188048573605
$ws.Range("F4").Value = 2143.181371125894
$ws.Range("G4").Value = 2127.424284502383
$ws.Range("H4").Value = 2205.988257890221
$ws.Range("B5").Value = 2180.197948453717
$ws.Range("C5").Value = 2025.331113879538
$ws.Range("D5").Value = 1344.466638394197
$ws.Range("E5").Value = 2212.0045224477
$ws.Range("F5").Value = 2184.578900396438
$ws.Range("G5").Value = 2120.356950798305
$ws.Range("H5").Value = 2212.657041173416
$ws.Range("B6").Value = 2186.955418566179
$ws.Range("C6").Value = 2024.724179085728
$ws.Range("D6").Value = 481.0406291909139
$ws.Range("E6").Value = 2216.909327032102
$ws.Range("F6").Value = 2190.875211451854
$ws.Range("G6").Value = 2056.101009838984
$ws.Range("H6").Value = 2217.150133649538
$ws.Range("B7").Value = 2166.472776123554
$ws.Range("C7").Value = 2030.493497846352
$ws.Range("D7").Value = 742.1278514679626
$ws.Range("E7").Value = 2209.28319458538
$ws.Range("F7").Value = 2172.310782951094
$ws.Range("G7").Value = 2070.712256698685
$ws.Range("H7").Value = 2210.270146677729
$ws.Range("B8").Value = 2167.740513980736
$ws.Range("C8").Value = 1982.384839167365
$ws.Range("D8").Value = 780.1099568704429
$ws.Range("E8").Value = 2202.028990745872
$ws.Range("F8").Value = 2169.76364612606
$ws.Range("G8").Value = 2049.817509468166
$ws.Range("H8").Value = 2202.406249820687
$ws.Range("B9").Value = 2181.940540741372
$ws.Range("C9").Value = 2019.818827124698
$ws.Range("D9").Value = 744.3932991661693
$ws.Range("E9").Value = 2205.775692409016
$ws.Range("F9").Value = 2189.436098052507
$ws.Range("G9").Value = 2073.482490100606
$ws.Range("H9").Value = 2207.082550820086
$ws.Range("B10").Value = 1908.84393146966
$ws.Range("C10").Value = 2038.806492436251
$ws.Range("D10").Value = 1529.336312426829
$ws.Range("E10").Value = 2172.380610578276
$ws.Range("F10").Value = 2024.861244482418
$ws.Range("G10").Value = 2114.98496527085
$ws.Range("H10").Value = 2178.023297273562
$ws.Range("B11").Value = 1850.571204774973
$ws.Range("C11").Value = 2047.386004418435
$ws.Range("D11").Value = 1373.274597388961
$ws.Range("E11").Value = 2173.445557392237
$ws.Range("F11").Value = 1967.172106832768
$ws.Range("G11").Value = 2103.257864896071
$ws.Range("H11").Value = 2177.081546657847
$ws.Range("B12").Value = 1557.822975848908
$ws.Range("C12").Value = 2028.861029587857
$ws.Range("D12").Value = 643.7507583034647
$ws.Range("E12").Value = 2134.639480504867
$ws.Range("F12").Value = 1609.080127237416
$ws.Range("G12").Value = 2043.878760777674
$ws.Range("H12").Value = 2135.293080965606
$ws.Range("B13").Value = 1889.092437598045
$ws.Range("C13").Value = 2034.465361463813
$ws.Range("D13").Value = 1332.77621821056
$ws.Range("E13").Value = 2167.399040879289
$ws.Range("F13").Value = 1979.265623244521
$ws.Range("G13").Value = 2086.80489728298
$ws.Range("H13").Value = 2170.737924912611
